$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells D1:H1 ---
$ws.Range("D1").Value = "pred_sd"
$ws.Range("E1").Value = "pred_lo_80"
$ws.Range("F1").Value = "pred_hi_80"
$ws.Range("G1").Value = "model"
$ws.Range("H1").Value = "exog_approval"

# Copy the bold/centered/bordered header style from A1 onto the new headers
$ws.Range("A1").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)

# --- Updated B/C values plus new D:H data values for rows 2-9 ---
$data = @(
    @{ Row=2; B=44.3382262630215;  C=1.782502371047188;  D=1.477564416905923;  E=42.44465127134272; F=46.23180125470029 },
    @{ Row=3; B=35.97098622709436; C=2.848716697536143;  D=2.635585571360043;  E=32.59334741199113; F=39.34862504219758 },
    @{ Row=4; B=2.723575158728821; C=0.6925965779035405; D=0.6317443829491085; E=1.913962155736383; F=3.533188161721258 },
    @{ Row=5; B=1.276709263833665; C=0.3026656212035301; D=0.2691008467461636; E=0.9318426523967411; F=1.621575875270588 },
    @{ Row=6; B=3.134954891566468; C=0.6201143687290629; D=0.6546931947237548; E=2.295931802916845; F=3.973977980216092 },
    @{ Row=7; B=2.049992320667111; C=0.7315696826819685; D=0.7002957758226334; E=1.152527172817345; F=2.947457468516877 },
    @{ Row=8; B=9.179588433454013; C=1.551212647681675;  D=1.478988215955746;  E=7.284188769873912; F=11.07498809703412 },
    @{ Row=9; B=1.15661508591859;  C=0.2959424388078578; D=0.2571016066550879; E=0.8271261194057304; F=1.48610405243145 }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = "ssm"
    $ws.Cells.Item($r, 8).Value = "on"
}
